$d = $word.ActiveDocument

# --- Replace the "proto file" paragraph text and relocate the _GoBack bookmark ---
$rng = $d.Content
$oldText = "The proto file must be explicitly named model.proto as Acumos expects it that way. Also note that package names must be globally unique to let AI4EU Experiments distinguish the protobuf definitions for all onboarded models.  "
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "  ", 2)
if (-not $found) {
    throw "Could not find the proto-file paragraph text to replace"
}
$bookmarkRange = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
Write-Host "Text + bookmark updated"

# --- Update OLE object ObjectID 1654001402 -> 1654349470 ---
$targetXml_1654001402 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00BC194B" w:rsidRPr="00F8750C" w:rsidRDefault="00AB0097" w:rsidP="00BC194B"><w:pPr><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:color w:val="002060"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:color w:val="002060"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:object w:dxaOrig="9360" w:dyaOrig="8649"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:468pt;height:432.75pt" o:ole=""><v:imagedata r:id="rId11" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1654349470" r:id="rId12"/></w:object></w:r></w:p>
'@
$found = $false
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*_1654001402*") {
        $p.Range.InsertXML($targetXml_1654001402)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find paragraph containing OLE object _1654001402"
}
Write-Host "Updated OLE object 1654001402"

# --- Update OLE object ObjectID 1654001403 -> 1654349471 ---
$targetXml_1654001403 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="002C1FF9" w:rsidRPr="002C1FF9" w:rsidRDefault="002C1FF9" w:rsidP="002C1FF9"><w:pPr><w:pStyle w:val="hj"/><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cstheme="minorHAnsi"/></w:rPr><w:object w:dxaOrig="9360" w:dyaOrig="9583"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:369.75pt;height:378.75pt" o:ole=""><v:imagedata r:id="rId13" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1026" DrawAspect="Content" ObjectID="_1654349471" r:id="rId14"/></w:object></w:r></w:p>
'@
$found = $false
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*_1654001403*") {
        $p.Range.InsertXML($targetXml_1654001403)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find paragraph containing OLE object _1654001403"
}
Write-Host "Updated OLE object 1654001403"

# --- Update OLE object ObjectID 1654001404 -> 1654349472 ---
$targetXml_1654001404 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00176708" w:rsidRPr="00F8750C" w:rsidRDefault="004F323E" w:rsidP="002C1FF9"><w:pPr><w:pStyle w:val="KeinLeerraum"/><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="00F8750C"><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cstheme="minorHAnsi"/></w:rPr><w:t>Below is the code snippet for client:</w:t></w:r><w:bookmarkStart w:id="3" w:name="_MON_1649246377"/><w:bookmarkEnd w:id="3"/><w:r w:rsidR="002C1FF9"><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cstheme="minorHAnsi"/></w:rPr><w:object w:dxaOrig="9360" w:dyaOrig="10285"><v:shape id="_x0000_i1027" type="#_x0000_t75" style="width:429.75pt;height:473.25pt" o:ole=""><v:imagedata r:id="rId16" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1027" DrawAspect="Content" ObjectID="_1654349472" r:id="rId17"/></w:object></w:r></w:p>
'@
$found = $false
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*_1654001404*") {
        $p.Range.InsertXML($targetXml_1654001404)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find paragraph containing OLE object _1654001404"
}
Write-Host "Updated OLE object 1654001404"

# --- Update OLE object ObjectID 1654001405 -> 1654349473 ---
$targetXml_1654001405 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00176708" w:rsidRDefault="002A7855" w:rsidP="00176708"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:color w:val="002060"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:color w:val="002060"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:object w:dxaOrig="9360" w:dyaOrig="5598"><v:shape id="_x0000_i1028" type="#_x0000_t75" style="width:468pt;height:279.75pt" o:ole=""><v:imagedata r:id="rId19" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1028" DrawAspect="Content" ObjectID="_1654349473" r:id="rId20"/></w:object></w:r></w:p>
'@
$found = $false
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*_1654001405*") {
        $p.Range.InsertXML($targetXml_1654001405)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find paragraph containing OLE object _1654001405"
}
Write-Host "Updated OLE object 1654001405"
